# ------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计",
#    populated with the quarters fund-holdings table.
# 2. Insert a new top data row in "总计" summarising the new quarter
#    and renumber its running index column.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$q4 = $wb.Worksheets.Item("2021-Q4")

# --- helper: write a value as genuine *text* (not auto-converted to a number) ---
# (kept as a convention in this script: set NumberFormat "@" before assigning the
#  value, then reset the cell style back to "Normal" once the text is in place so
#  no left-over number-format shows up on the cell)

# ================= Step 1: new sheet "2022-Q1" =================
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4)
$q1.Name = "2022-Q1"

# NOTE: fetch the "总计" sheet reference only *after* inserting the new sheet -
# worksheet handles are resolved positionally, so grabbing it beforehand would
# silently end up pointing at the freshly inserted sheet instead.
$zj = $wb.Worksheets.Item("总计")

# -- header row --
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# copy the bold/centered header style (and the index-column style) used by the
# other quarterly sheets so the new sheet matches their look
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# -- data rows --
$q1.Range("A2").Value = 0
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "010490"
$q1.Range("C2").Value = "鹏华高质量增长混合A"
$q1.Range("D2").Value = "13.31"
$q1.Range("E2").Value = "93.61"
$q1.Range("F2").Value = "8.54"
$q1.Range("G2").Value = "1.1367"
$q1.Range("B2:G2").Style = "Normal"
$q1.Range("H2").Value = 2

$q1.Range("A3").Value = 1
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "160607"
$q1.Range("C3").Value = "鹏华价值优势混合(LOF)"
$q1.Range("D3").Value = "15.59"
$q1.Range("E3").Value = "89.56"
$q1.Range("F3").Value = "5.25"
$q1.Range("G3").Value = "0.8185"
$q1.Range("B3:G3").Style = "Normal"
$q1.Range("H3").Value = 2

$q1.Range("A4").Value = 2
$q1.Range("B4:G4").NumberFormat = "@"
$q1.Range("B4").Value = "009023"
$q1.Range("C4").Value = "鹏华稳健回报混合"
$q1.Range("D4").Value = "3.52"
$q1.Range("E4").Value = "93.91"
$q1.Range("F4").Value = "10.07"
$q1.Range("G4").Value = "0.3545"
$q1.Range("B4:G4").Style = "Normal"
$q1.Range("H4").Value = 1

$q1.Range("A5").Value = 3
$q1.Range("B5:G5").NumberFormat = "@"
$q1.Range("B5").Value = "501062"
$q1.Range("C5").Value = "南方瑞合三年定期开放混合(LOF)"
$q1.Range("D5").Value = "6.88"
$q1.Range("E5").Value = "56.71"
$q1.Range("F5").Value = "2.81"
$q1.Range("G5").Value = "0.1933"
$q1.Range("B5:G5").Style = "Normal"
$q1.Range("H5").Value = 10

$q1.Range("A6").Value = 4
$q1.Range("B6:G6").NumberFormat = "@"
$q1.Range("B6").Value = "006976"
$q1.Range("C6").Value = "鹏华核心优势混合"
$q1.Range("D6").Value = "2.40"
$q1.Range("E6").Value = "91.56"
$q1.Range("F6").Value = "7.60"
$q1.Range("G6").Value = "0.1824"
$q1.Range("B6:G6").Style = "Normal"
$q1.Range("H6").Value = 1

$q1.Range("A7").Value = 5
$q1.Range("B7:G7").NumberFormat = "@"
$q1.Range("B7").Value = "008811"
$q1.Range("C7").Value = "鹏华科技创新混合"
$q1.Range("D7").Value = "3.49"
$q1.Range("E7").Value = "94.05"
$q1.Range("F7").Value = "2.85"
$q1.Range("G7").Value = "0.0995"
$q1.Range("B7:G7").Style = "Normal"
$q1.Range("H7").Value = 10

$q1.Range("A8").Value = 6
$q1.Range("B8:G8").NumberFormat = "@"
$q1.Range("B8").Value = "010491"
$q1.Range("C8").Value = "鹏华高质量增长混合C"
$q1.Range("D8").Value = "0.28"
$q1.Range("E8").Value = "93.61"
$q1.Range("F8").Value = "8.54"
$q1.Range("G8").Value = "0.0239"
$q1.Range("B8:G8").Style = "Normal"
$q1.Range("H8").Value = 2

$q1.Range("A9").Value = 7
$q1.Range("B9:G9").NumberFormat = "@"
$q1.Range("B9").Value = "003835"
$q1.Range("C9").Value = "鹏华沪深港新兴成长灵活配置混合"
$q1.Range("D9").Value = "0.61"
$q1.Range("E9").Value = "82.70"
$q1.Range("F9").Value = "3.54"
$q1.Range("G9").Value = "0.0216"
$q1.Range("B9:G9").Style = "Normal"
$q1.Range("H9").Value = 7

# ================= Step 2: update the "总计" summary sheet =================
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()

# give the new index cell (A2) the same style as the other index cells
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 8
$zj.Range("D2").Value = 2.83

# renumber the running index column (0, 1, 2, ...) for the rows pushed down
$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2

Write-Host "Added 2022-Q1 sheet and updated 总计 summary"
